# Commit: Fri, Jun 12, 2020  8:04:59 PM
#
# Two things happened in this commit:
#
# 1. The deck's three tables (on slides 14, 15 and 16) had their table
#    style switched from the locally-defined "Table_0" style
#    ({85E2D9D4-BA40-4531-8218-082E24BA75A2}) to the built-in table
#    style {9A17EA7A-D257-4DAD-A51A-B0156E223351}.
#
# 2. The presentation's colour theme was changed from the colourful
#    "Integral" / "Red Violet" scheme back to the plain default
#    "Office Theme" colour scheme (the deck's font scheme and format
#    scheme were already the stock "Office" ones, so only the colours
#    move).

$p = $ppt.ActivePresentation

# --- 1. Table styles -----------------------------------------------------

$oldStyleId = "{85E2D9D4-BA40-4531-8218-082E24BA75A2}"
$newStyleId = "{9A17EA7A-D257-4DAD-A51A-B0156E223351}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)

        if ($shape.HasTable) {
            $table = $shape.Table

            if ($table.Style.Name -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Theme colours ------------------------------------------------------
# ThemeColorScheme items are ordered dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink - switch every slot back to the standard "Office Theme" palette.
# (RGB is encoded VBA-style as R + G*256 + B*65536.)

$officeThemeColors = 0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($colorIdx = 1; $colorIdx -le $themeColors.Count; $colorIdx++) {
    $themeColors.Item($colorIdx).RGB = $officeThemeColors[$colorIdx - 1]
}
